# Work log update: add two new work items, fix a typo in an existing item.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new row before row 89 for the "Java utility classes" task ---
$ws.Rows.Item(89).Insert()

$ws.Range("A89").Value = "Java utility classes in separate project to manipulate the game board tile positions, scaling in JSON "
$ws.Range("B89").Value = "UI"
$ws.Range("C89").Value = [DateTime]"2018-12-12"
$ws.Range("D89").Value = "Andrew"
$ws.Range("E89").Value = [DateTime]"2018-12-13"
$ws.Range("F89").Value = "Board displays correctly in UI"

# --- Fix a typo in the item that got pushed down to row 90 ---
$ws.Range("A90").Value = "CareerChange prompts a player with a CollegeCareer with StandardCareer options"

# --- Append a new row at the end (row 93) for the "GameBoard refactor" task ---
$ws.Range("A93").Value = "Refactoring of GameBoard initialisation to use façade design pattern"
$ws.Range("B93").Value = "Refactor"
$ws.Range("C93").Value = [DateTime]"2018-12-19"
$ws.Range("D93").Value = "Andrew"
$ws.Range("E93").Value = [DateTime]"2018-12-19"
$ws.Range("F93").Value = "Completed, unit and integration tests all passing"

# --- Restore the date-format style on the new row's date cells (matches existing rows) ---
$ws.Range("C89").NumberFormat = $ws.Range("C92").NumberFormat
$ws.Range("E89").NumberFormat = $ws.Range("E92").NumberFormat
$ws.Range("C93").NumberFormat = $ws.Range("C92").NumberFormat
$ws.Range("E93").NumberFormat = $ws.Range("E92").NumberFormat

# --- Mirror the cursor/selection move to the newly active cell ---
$ws.Range("A91").Select()

Write-Host "Work log updated."
